$d = $word.ActiveDocument

$d.Content.Find.Execute("924÷2=462, 0", $true, $false, $false, $false, $false, $true, 1, $false, "825÷4=206, 1", 2) | Out-Null
$d.Content.Find.Execute("889÷2=444, 1", $true, $false, $false, $false, $false, $true, 1, $false, "469÷7=67, 0", 2) | Out-Null
$d.Content.Find.Execute("370÷2=185, 0", $true, $false, $false, $false, $false, $true, 1, $false, "291÷5=58, 1", 2) | Out-Null
$d.Content.Find.Execute("977÷7=139, 4", $true, $false, $false, $false, $false, $true, 1, $false, "656÷8=82, 0", 2) | Out-Null
$d.Content.Find.Execute("470÷6=78, 2", $true, $false, $false, $false, $false, $true, 1, $false, "776÷7=110, 6", 2) | Out-Null
$d.Content.Find.Execute("460÷3=153, 1", $true, $false, $false, $false, $false, $true, 1, $false, "252÷9=28, 0", 2) | Out-Null
$d.Content.Find.Execute("111÷8=13, 7", $true, $false, $false, $false, $false, $true, 1, $false, "730÷8=91, 2", 2) | Out-Null
$d.Content.Find.Execute("836÷2=418, 0", $true, $false, $false, $false, $false, $true, 1, $false, "824÷4=206, 0", 2) | Out-Null
$d.Content.Find.Execute("409÷8=51, 1", $true, $false, $false, $false, $false, $true, 1, $false, "185÷3=61, 2", 2) | Out-Null
$d.Content.Find.Execute("452÷8=56, 4", $true, $false, $false, $false, $false, $true, 1, $false, "630÷7=90, 0", 2) | Out-Null
$d.Content.Find.Execute("351÷4=87, 3", $true, $false, $false, $false, $false, $true, 1, $false, "960÷4=240, 0", 2) | Out-Null
$d.Content.Find.Execute("741÷6=123, 3", $true, $false, $false, $false, $false, $true, 1, $false, "696÷6=116, 0", 2) | Out-Null
$d.Content.Find.Execute("381÷9=42, 3", $true, $false, $false, $false, $false, $true, 1, $false, "374÷7=53, 3", 2) | Out-Null
$d.Content.Find.Execute("296÷9=32, 8", $true, $false, $false, $false, $false, $true, 1, $false, "357÷6=59, 3", 2) | Out-Null
$d.Content.Find.Execute("669÷2=334, 1", $true, $false, $false, $false, $false, $true, 1, $false, "878÷6=146, 2", 2) | Out-Null
$d.Content.Find.Execute("531÷2=265, 1", $true, $false, $false, $false, $false, $true, 1, $false, "551÷6=91, 5", 2) | Out-Null
$d.Content.Find.Execute("335÷2=167, 1", $true, $false, $false, $false, $false, $true, 1, $false, "883÷7=126, 1", 2) | Out-Null
$d.Content.Find.Execute("815÷8=101, 7", $true, $false, $false, $false, $false, $true, 1, $false, "573÷2=286, 1", 2) | Out-Null
$d.Content.Find.Execute("486÷2=243, 0", $true, $false, $false, $false, $false, $true, 1, $false, "412÷3=137, 1", 2) | Out-Null
$d.Content.Find.Execute("865÷8=108, 1", $true, $false, $false, $false, $false, $true, 1, $false, "819÷3=273, 0", 2) | Out-Null
$d.Content.Find.Execute("381÷6=63, 3", $true, $false, $false, $false, $false, $true, 1, $false, "354÷2=177, 0", 2) | Out-Null
$d.Content.Find.Execute("507÷2=253, 1", $true, $false, $false, $false, $false, $true, 1, $false, "442÷2=221, 0", 2) | Out-Null
$d.Content.Find.Execute("679÷8=84, 7", $true, $false, $false, $false, $false, $true, 1, $false, "335÷8=41, 7", 2) | Out-Null
$d.Content.Find.Execute("472÷5=94, 2", $true, $false, $false, $false, $false, $true, 1, $false, "524÷3=174, 2", 2) | Out-Null
$d.Content.Find.Execute("676÷3=225, 1", $true, $false, $false, $false, $false, $true, 1, $false, "164÷6=27, 2", 2) | Out-Null
